$ErrorActionPreference = "Stop"
$d = $word.ActiveDocument

function Get-ParagraphByText($doc, $text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        $t = $p.Range.Text
        $t = $t.TrimEnd([char]7, [char]13, [char]12)
        if ($t -eq $text) {
            return $p
        }
    }
    return $null
}

# --- Edit 1: after the "Net cham gach" heading, insert the two new
#     descriptive paragraphs (algorithm used / point changed) ---
$headingNetChamGach = Get-ParagraphByText $d "Nét chấm gạch"
$afterHeadingRange = $headingNetChamGach.Range.InsertParagraphAfter()
$newPara1 = Get-ParagraphByText $d "Nét chấm gạch"
$newPara1 = $newPara1.Next()
$xml1 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Thuật toán sử dụng:</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> Thuật toán vẽ đường thẳng </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Midpoint</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Điểm thay đổi:</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Tương tự đường đứt nét. Đường chấm gạch cũng sử dụng biến đếm để loại bỏ các điểm không cần thiết.VD: tại vị trí x-1 thì bỏ qua không putpixel kế tiếp đến x thì được put , x+1 không được put vậy ta có được nét chấm , còn các đoạn liền nét nằm ngoài điều kiện của biến đếm thì được putpixel bình thường </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$newPara1.Range.InsertXML($xml1)

# --- Edit 2: after the "Net 2 cham gach" heading, insert a blank bold
#     paragraph (spacer) ---
$headingNet2ChamGach = Get-ParagraphByText $d "Nét 2 chấm gạch"
$headingNet2ChamGach.Range.InsertParagraphAfter()
$newPara2 = Get-ParagraphByText $d "Nét 2 chấm gạch"
$newPara2 = $newPara2.Next()
$xml2 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$newPara2.Range.InsertXML($xml2)

# --- Edit 3: mark the second inline drawing (Picture 2) as the start
#     of a rendered page, i.e. add <w:lastRenderedPageBreak/> right
#     before its <w:drawing> element ---
$shape2 = $d.InlineShapes.Item(2)
$drawRange = $shape2.Range
$drawRange.Collapse(1)
$lprXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p><w:r><w:rPr><w:noProof/></w:rPr><w:lastRenderedPageBreak/></w:r></w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@
$drawRange.InsertXML($lprXml)

Write-Output "done"
